$d = $word.ActiveDocument

# 1. "Defendant appeared in Court on July 03, 2022" -> "... July 04, 2022"
$d.Content.Find.Execute("Defendant appeared in Court on July 03, 2022", $false, $false, $false, $false, $false,
                         $true, 0, $false, "Defendant appeared in Court on July 04, 2022", 1)

# 2. standalone "July 03, 2022" -> "July 04, 2022" (the bold date line)
$d.Content.Find.Execute("July 03, 2022", $false, $false, $false, $false, $false,
                         $true, 0, $false, "July 04, 2022", 1)

# 3. "September 01, 2022" -> "September 02, 2022"
$d.Content.Find.Execute("September 01, 2022", $false, $false, $false, $false, $false,
                         $true, 0, $false, "September 02, 2022", 1)

# 4. " license is suspended from July 03, 2022" -> "... July 04, 2022"
$d.Content.Find.Execute(" license is suspended from July 03, 2022", $false, $false, $false, $false, $false,
                         $true, 0, $false, " license is suspended from July 04, 2022", 1)
